$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.550.14"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "3.391.03"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'575.66"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "'141.15"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").Value = "'7.66"
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -2.68%  "
$ws.Range("D12").Value = "3.969.05"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "3.396.60"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "61.597.75"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'13.61"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "'8.99"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'391.22"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "'75.00"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'0.0000112"
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("E26").Value = "  +7.59%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "'8.03"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.39"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'6.90"
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'168.69"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "3.424.93"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -1.43%  "
$ws.Range("D39").Value = "'0.0765"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("E40").Value = "  -3.92%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'4.41"
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").Value = "2.482.06"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("D46").Value = "'22.85"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("D47").Value = "'6.66"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'2.03"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("D51").Value = "'0.206"
$ws.Range("E51").Value = "  -2.13%  "
